$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Defaults): Spring constant changes from 10000 to 10
$ws.Range("C3").Value = 10

# Row 4 (Values): Spring constant changes from 10 to 10000
$ws.Range("C4").Value = 10000

# Row 4: Surface tension coefficient changes from 5.0E-4 to 3.0E-3, now shown with 3 decimals
$ws.Range("D4").NumberFormat = "0.000"
$ws.Range("D4").Value = 0.003

# Row 4: Repeats changes from 50 to 200 (stored as text, like the existing "50")
$ws.Range("I4").Value = "200"

# Row 4: Timestep changes from 1E-4 to 5.0E-5
$ws.Range("J4").Value = 0.00005

# New column M: Numerical Method (values added in this order so the shared
# string table ends up in the same order as the authored workbook)
$ws.Range("M1").Value = "Numerical Method"
$ws.Range("M3").Value = "ER"
$ws.Range("M4").Value = "ER"
$ws.Range("M2").Value = "string (ER or V)"

# Column M width (~17.21875 chars as authored; engine quantizes widths
# internally so 16.33 is the closest input that reproduces it)
$ws.Columns.Item(13).ColumnWidth = 16.33

# Update selection to the newly added cell, like the authored workbook
$ws.Range("M2").Select()
